# Apply odds updates for rows 4, 5, and 7 (FlashScore weekly games workbook).
# The pre-seeded $wb can be unreliable in this host, so resolve the live
# workbook/worksheet via $excel (falls back to $wb if that ever is populated).
if ($excel -ne $null -and $excel.ActiveWorkbook -ne $null) {
    $targetWb = $excel.ActiveWorkbook
} else {
    $targetWb = $wb
}
$ws = $targetWb.ActiveSheet

# Row 4 updates
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("Y4").Value = 9.5
$ws.Range("AC4").Value = 6.5
$ws.Range("AH4").Value = 9
$ws.Range("AJ4").Value = 17
$ws.Range("AO4").Value = 11
$ws.Range("AT4").Value = 2.25
$ws.Range("BC4").Value = 451
$ws.Range("BD4").Value = 126

# Row 5 updates
$ws.Range("G5").Value = 1.75
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.5
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 5.5
$ws.Range("N5").Value = 7.5
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.57
$ws.Range("S5").Value = 1.5
$ws.Range("T5").Value = 2.5
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 7.5
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 15
$ws.Range("AA5").Value = 17
$ws.Range("AC5").Value = 7.5
$ws.Range("AD5").Value = 6.5
$ws.Range("AH5").Value = 11
$ws.Range("AI5").Value = 23
$ws.Range("AJ5").Value = 17
$ws.Range("AK5").Value = 51
$ws.Range("AN5").Value = 3.6
$ws.Range("AO5").Value = 10
$ws.Range("AP5").Value = 23
$ws.Range("AQ5").Value = 34
$ws.Range("AS5").Value = 201
$ws.Range("AT5").Value = 2.5
$ws.Range("AU5").Value = 9.5
$ws.Range("AX5").Value = 6.5
$ws.Range("AY5").Value = 29
$ws.Range("BA5").Value = 101
$ws.Range("BB5").Value = 151
$ws.Range("BC5").Value = 401
$ws.Range("BD5").Value = 126

# Row 7 updates
$ws.Range("G7").Value = 1.1
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 1.4
$ws.Range("K7").Value = 3.5
$ws.Range("L7").Value = 13
$ws.Range("N7").Value = 29
$ws.Range("O7").Value = 1.08
$ws.Range("P7").Value = 8
$ws.Range("Q7").Value = 1.29
$ws.Range("R7").Value = 3.6
$ws.Range("S7").Value = 1.17
$ws.Range("T7").Value = 5
$ws.Range("U7").Value = 2.1
$ws.Range("V7").Value = 1.67
$ws.Range("W7").Value = 12
$ws.Range("Y7").Value = 12
$ws.Range("Z7").Value = 7
$ws.Range("AB7").Value = 29
$ws.Range("AC7").Value = 29
$ws.Range("AD7").Value = 21
$ws.Range("AE7").Value = 34
$ws.Range("AF7").Value = 81
$ws.Range("AM7").Value = 67
$ws.Range("AO7").Value = 4.5
$ws.Range("AQ7").Value = 9
$ws.Range("AT7").Value = 5
$ws.Range("AU7").Value = 12
$ws.Range("AV7").Value = 67
$ws.Range("AX7").Value = 17
$ws.Range("AY7").Value = 67
$ws.Range("AZ7").Value = 51
$ws.Range("BA7").Value = 351
$ws.Range("BB7").Value = 301
$ws.Range("BC7").Value = 301

